$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# --- survey sheet: add a new "select_multiple reason" question row ---
$survey.Range("A3").Value = "select_multiple reason"
$survey.Range("B3").Value = "visit_reason"
$survey.Range("C3").Value = "Reason of visits"

# --- choices sheet: add the "reason" choice list with 4 options ---
$choices.Range("A2").Value = "reason"
$choices.Range("B2").Value = "pregnant"
$choices.Range("C2").Value = "Pregnant"

$choices.Range("A3").Value = "reason"
$choices.Range("B3").Value = "checkup"
$choices.Range("C3").Value = "Check-ups"

$choices.Range("A4").Value = "reason"
$choices.Range("B4").Value = "vaccination"
$choices.Range("C4").Value = "Vaccination"

$choices.Range("A5").Value = "reason"
$choices.Range("B5").Value = "malnutrition"
$choices.Range("C5").Value = "Malnutrition"

# --- selection / active sheet bookkeeping ---
$survey.Activate() | Out-Null
$survey.Range("C4").Select() | Out-Null

$choices.Activate() | Out-Null
$choices.Range("C6").Select() | Out-Null
